# Append a new data row (row 72) to the active sheet, mirroring the
# existing Adafruit IO export rows (timestamp, feed key, value, lat/long/elev).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 72

$ws.Range("A$newRow").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$newRow").Value = "temperature"

# Keep the numeric-looking "Value" column stored as text, matching the
# existing rows (e.g. C71 = "25" as text, not a number).
$ws.Range("C$newRow").NumberFormat = "@"
$ws.Range("C$newRow").Value = "25"

$ws.Range("D$newRow").Value = "N/A"
$ws.Range("E$newRow").Value = "N/A"
$ws.Range("F$newRow").Value = "N/A"
